# Updated symbol list on Tue Jan 24 14:39:44 UTC 2023 with GitHub Actions
# This script applies updated Price (D) and Volume(1h) (E) text values
# to the cryptocurrency listing on Sheet1, matching the source diff exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (kept as literal text, not
# auto-converted to a number/percentage, matching the original inline-string cells)
$updates = @{
    "D2" = "314.32"
    "E2" = "3.31%"
    "D3" = "35.75"
    "E3" = "0.71%"
    "D4" = "5.092"
    "E4" = "0.61%"
    "D5" = "0.08136"
    "E5" = "3.13%"
    "D6" = "2.103"
    "E6" = "-0.08%"
    "D7" = "4.140"
    "E7" = "0.31%"
    "D8" = "7.945"
    "E8" = "0.53%"
    "D9" = "0.9324"
    "E9" = "1.09%"
    "D10" = "0.1030"
    "E10" = "6.26%"
    "D11" = "0.1925"
    "E11" = "4.65%"
    "D12" = "0.09008"
    "E12" = "4.09%"
    "D13" = "0.03607"
    "E13" = "1.37%"
    "D14" = "0.09889"
    "E14" = "-0.40%"
    "D15" = "0.001444"
    "E15" = "0.93%"
    "D16" = "0.005805"
    "E16" = "2.69%"
    "D17" = "3.465"
    "E17" = "-0.19%"
    "D18" = "2.887"
    "E18" = "9.35%"
    "E19" = "1.12%"
    "D20" = "0.1312"
    "E20" = "-2.03%"
    "D21" = "5.106"
    "E22" = "0.10%"
    "D23" = "0.04555"
    "E23" = "1.11%"
    "D24" = "0.001248"
    "E24" = "1.09%"
    "D25" = "0.004794"
    "E25" = "-1.19%"
    "D26" = "0.0001253"
    "E26" = "-3.93%"
    "D27" = "0.0004512"
    "E27" = "-5.31%"
    "D39" = "0.01959"
    "E39" = "6.92%"
    "D40" = "0.04881"
    "E40" = "3.43%"
    "D41" = "0.007626"
    "E41" = "-3.44%"
    "D42" = "0.1382"
    "E42" = "-0.53%"
    "D43" = "0.007763"
    "E43" = "0.09%"
    "D44" = "0.002108"
    "E44" = "-4.00%"
    "D45" = "0.01175"
    "E45" = "5.30%"
    "D46" = "0.00006718"
    "E46" = "5.40%"
    "E47" = "-0.08%"
    "D48" = "156.80"
    "E48" = "210.61%"
    "D49" = "0.001705"
    "E49" = "-10.59%"
    "E50" = "-0.08%"
    "E51" = "-0.08%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"   # force text storage so numeric-looking strings are not reinterpreted
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"     # drop the temporary text format, restoring the original (no explicit) style
}
